# Auto-generated edit script: updates currentAveragePrice/LevePrice/LeveProfit
# columns (H-N) on each class sheet to refreshed market-board figures, as
# produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 247.14285
$ws.Range("I6").Value = 188.16667
$ws.Range("J6").Value = 601
$ws.Range("K6").Value = 564.50001
$ws.Range("L6").Value = 1803
$ws.Range("M6").Value = -452.50001
$ws.Range("N6").Value = -2027
$ws.Range("H12").Value = 102.5
$ws.Range("I12").Value = 103.333336
$ws.Range("K12").Value = 103.333336
$ws.Range("M12").Value = 66.666664
$ws.Range("H32").Value = 2052.6667
$ws.Range("I32").Value = 996
$ws.Range("K32").Value = 996
$ws.Range("M32").Value = -670
$ws.Range("H38").Value = 1758.909
$ws.Range("I38").Value = 415
$ws.Range("J38").Value = 5342.6665
$ws.Range("K38").Value = 1245
$ws.Range("L38").Value = 16027.9995
$ws.Range("M38").Value = -873
$ws.Range("N38").Value = -16771.9995
$ws.Range("H55").Value = 804.06665
$ws.Range("J55").Value = 942.2857
$ws.Range("L55").Value = 942.2857
$ws.Range("N55").Value = -1370.2857
$ws.Range("H111").Value = 1946.4
$ws.Range("I111").Value = 1946.4
$ws.Range("K111").Value = 5839.200000000001
$ws.Range("M111").Value = -2772.200000000001
$ws.Range("H132").Value = 2204.2
$ws.Range("I132").Value = 2242.75
$ws.Range("K132").Value = 6728.25
$ws.Range("M132").Value = -4198.25
$ws.Range("H138").Value = 4973.289
$ws.Range("I138").Value = 2807.2646
$ws.Range("K138").Value = 8421.793799999999
$ws.Range("M138").Value = -3281.793799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2297.8333
$ws.Range("I2").Value = 998.5
$ws.Range("K2").Value = 998.5
$ws.Range("M2").Value = -885.5
$ws.Range("H32").Value = 7059.175
$ws.Range("I32").Value = 4781.5586
$ws.Range("K32").Value = 4781.5586
$ws.Range("M32").Value = -4494.5586
$ws.Range("H45").Value = 2541
$ws.Range("I45").Value = 2648.4546
$ws.Range("J45").Value = 1950
$ws.Range("K45").Value = 2648.4546
$ws.Range("L45").Value = 1950
$ws.Range("M45").Value = -2271.4546
$ws.Range("N45").Value = -2704
$ws.Range("H95").Value = 31596
$ws.Range("J95").Value = 31596
$ws.Range("L95").Value = 31596
$ws.Range("N95").Value = -37088
$ws.Range("H102").Value = 907.5
$ws.Range("I102").Value = 907.5
$ws.Range("K102").Value = 907.5
$ws.Range("M102").Value = 714.5
$ws.Range("H110").Value = 3735.5386
$ws.Range("I110").Value = 2667.2
$ws.Range("K110").Value = 2667.2
$ws.Range("M110").Value = -622.1999999999998
$ws.Range("H116").Value = 2297.8333
$ws.Range("I116").Value = 998.5
$ws.Range("K116").Value = 998.5
$ws.Range("M116").Value = 1295.5
$ws.Range("H122").Value = 3410.4443
$ws.Range("I122").Value = 3738.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11216.4
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8766.400000000001
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 1603
$ws.Range("I132").Value = 1629.5
$ws.Range("K132").Value = 4888.5
$ws.Range("M132").Value = -2358.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2297.8333
$ws.Range("I3").Value = 998.5
$ws.Range("K3").Value = 998.5
$ws.Range("M3").Value = -884.5
$ws.Range("H22").Value = 421.1111
$ws.Range("I22").Value = 357.85715
$ws.Range("J22").Value = 642.5
$ws.Range("K22").Value = 357.85715
$ws.Range("L22").Value = 642.5
$ws.Range("M22").Value = -184.85715
$ws.Range("N22").Value = -988.5
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
$ws.Range("H99").Value = 3300
$ws.Range("I99").Value = 3300
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3300
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1802
$ws.Range("N99").Value = $null
$ws.Range("H105").Value = 2503.1667
$ws.Range("I105").Value = 2325.5715
$ws.Range("J105").Value = 3124.75
$ws.Range("K105").Value = 2325.5715
$ws.Range("L105").Value = 3124.75
$ws.Range("M105").Value = -578.5715
$ws.Range("N105").Value = -6618.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 95.78570999999999
$ws.Range("I7").Value = 45.3
$ws.Range("K7").Value = 45.3
$ws.Range("M7").Value = 67.7
$ws.Range("H92").Value = 33875
$ws.Range("J92").Value = 33875
$ws.Range("L92").Value = 33875
$ws.Range("N92").Value = -38867
$ws.Range("H99").Value = 13720.863
$ws.Range("J99").Value = 14368
$ws.Range("L99").Value = 14368
$ws.Range("N99").Value = -17364
$ws.Range("H126").Value = 13720.863
$ws.Range("J126").Value = 14368
$ws.Range("L126").Value = 43104
$ws.Range("N126").Value = -48044

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 8569.333000000001
$ws.Range("I97").Value = 7995
$ws.Range("J97").Value = 8856.5
$ws.Range("K97").Value = 23985
$ws.Range("L97").Value = 26569.5
$ws.Range("M97").Value = -23489
$ws.Range("N97").Value = -27561.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5909.4546
$ws.Range("I70").Value = 5002.143
$ws.Range("K70").Value = 5002.143
$ws.Range("M70").Value = -4732.143
$ws.Range("H73").Value = 5909.4546
$ws.Range("I73").Value = 5002.143
$ws.Range("K73").Value = 5002.143
$ws.Range("M73").Value = -4066.143
$ws.Range("H122").Value = 61271.94
$ws.Range("I122").Value = 1638
$ws.Range("K122").Value = 4914
$ws.Range("M122").Value = -2464
$ws.Range("H132").Value = 2742.3333
$ws.Range("I132").Value = 2213.3333
$ws.Range("K132").Value = 6639.999899999999
$ws.Range("M132").Value = -4109.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2448.875
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
$ws.Range("H26").Value = 5924.5
$ws.Range("I26").Value = 5924.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 5924.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -5629.5
$ws.Range("N26").Value = $null
$ws.Range("H61").Value = 4181
$ws.Range("I61").Value = 4387.8887
$ws.Range("J61").Value = 3250
$ws.Range("K61").Value = 4387.8887
$ws.Range("L61").Value = 3250
$ws.Range("M61").Value = -4185.8887
$ws.Range("N61").Value = -3654
$ws.Range("H113").Value = 4181
$ws.Range("I113").Value = 4387.8887
$ws.Range("J113").Value = 3250
$ws.Range("K113").Value = 4387.8887
$ws.Range("L113").Value = 3250
$ws.Range("M113").Value = -2217.8887
$ws.Range("N113").Value = -7590
$ws.Range("H126").Value = 2448.875
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1924.625
$ws.Range("I122").Value = 1876.5555
$ws.Range("J122").Value = 1986.4286
$ws.Range("K122").Value = 5629.666499999999
$ws.Range("L122").Value = 5959.2858
$ws.Range("M122").Value = -3179.666499999999
$ws.Range("N122").Value = -10859.2858
